$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metricas")

# Row 10: new task "implementar PilaHL"
$ws.Range("A10").Value = "implementar PilaHL"
$ws.Range("B10").Value = 20
$ws.Range("D10").Value = 0.0069444444444444441
$ws.Range("E10").Value = 0.41666666666666669
$ws.Range("F10").Value = 0.42222222222222222
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 0.00069444444444444447

# Move the active cell selection to A14 (as captured in the saved view state)
$ws.Range("A14").Select()

$excel.Calculate()
